# "Fixed typo in Tiers table"
# - Correct "refraction-limited" -> "diffraction-limited" in the D4 cell
#   (Tier 2 / Advanced Quantification description) on all three sheets.
# - Restore the print area of the main "Tier system_v02-00" sheet to
#   $A$1:$L$5 (it had drifted to $A$2:$I$5).
# - Re-activate the "Tier system_v02-00" sheet/tab and update the saved
#   cell selections on each sheet to match the author's last editing
#   position.

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "Tier system_v02-00",
    "Tier system_v02-00 SUMMARY",
    "Tier system_v02-00 MINIMAL"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $cell = $ws.Range("D4")
    $text = $cell.Value()
    $fixed = $text.Replace("refraction-limited", "diffraction-limited")
    $cell.Value = $fixed
}

# Restore the print area on the first (main) sheet.
$ws1 = $wb.Worksheets.Item("Tier system_v02-00")
$ws1.PageSetup.PrintArea = '$A$1:$L$5'

# Update view state: make "Tier system_v02-00" the active/selected tab
# again, and set each sheet's remembered selection.
$ws1.Activate()
$ws1.Range("H4").Select()

$ws2 = $wb.Worksheets.Item("Tier system_v02-00 SUMMARY")
$ws2.Activate()
$ws2.Range("G4").Select()

$ws3 = $wb.Worksheets.Item("Tier system_v02-00 MINIMAL")
$ws3.Activate()
$ws3.Range("E8").Select()

$ws1.Activate()
